# Regenerate the s_val data (filters save games), updating columns B:E and G
# for rows 2-11. Column A (dates) and F (Win) remain unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(1.455362044514542,  1.655778082260271,  0.1494219747398047, 0.4942365360607697, 3.754798637575387)
    3  = @(3.286832544864788,  1.655778082260271,  3.537761648806719,  0.4942365360607697, 8.974608811992548)
    4  = @(0.1190320826869504, 0.306821227259698,  22.3905356188092,   10.19245300693656,  33.0088419356924)
    5  = @(3.286832544864788,  1.655778082260271,  0.7527432677738641, 0.4942365360607697, 6.189590430959694)
    6  = @(0.1190320826869504, 0.04071648406533734,3.537761648806719,  0.4942365360607697, 4.191746751619776)
    7  = @(0.1190320826869504, 0.002571899574220771,0.7527432677738641,0.4942365360607697, 1.368583786095805)
    8  = @(0.1190320826869504, 0.306821227259698,  0.7527432677738641, 10.19245300693656,  11.37104958465707)
    9  = @(3.286832544864788,  1.655778082260271,  0.1494219747398047, 0.4942365360607697, 5.586269137925634)
    10 = @(0.1190320826869504, 0.04071648406533734,0.1494219747398047, 0.4942365360607697, 0.8034070775528621)
    11 = @(3.286832544864788,  1.655778082260271,  0.7527432677738641, 0.4942365360607697, 6.189590430959694)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]  # B - TB
    $ws.Cells.Item($row, 3).Value = $vals[1]  # C - d2S
    $ws.Cells.Item($row, 4).Value = $vals[2]  # D - K
    $ws.Cells.Item($row, 5).Value = $vals[3]  # E - IP
    $ws.Cells.Item($row, 7).Value = $vals[4]  # G - sum
}
